# Update database: shift quarterly/annual periods one slot to the left
# (drop oldest "6 ماهه منتهی به 1399/06" column, append newest
# "12 ماهه منتهی به 1401/12"), and correct a couple of historical
# "read_price" values (capital + the two EPS rows that derive from it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","E","F","G","H","I","J","K","L","M")

# ---------------------------------------------------------------
# Row 8: فترة مالی (financial period) headers
# ---------------------------------------------------------------
$periods = @(
    "9 ماهه منتهی به 1399/09",
    "12 ماهه منتهی به 1399/12",
    "3 ماهه منتهی به 1400/03",
    "6 ماهه منتهی به 1400/06",
    "9 ماهه منتهی به 1400/09",
    "12 ماهه منتهی به 1400/12",
    "3 ماهه منتهی به 1401/03",
    "6 ماهه منتهی به 1401/06",
    "9 ماهه منتهی به 1401/09",
    "12 ماهه منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periods[$i]
}

# ---------------------------------------------------------------
# Row 9: تاریخ انتشار (publish date) headers
# ---------------------------------------------------------------
$dates = @(
    "1400-10-29 (3)",
    "1401-02-10 (12)",
    "1401-05-12 (4)",
    "1401-09-09 (4)",
    "1401-10-29 (3)",
    "1402-02-09 (10)",
    "1401-05-12 (2)",
    "1401-09-09 (2)",
    "1401-10-29",
    "1402-02-09 (2)"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "9").Value = $dates[$i]
}

# ---------------------------------------------------------------
# Data rows 11-27 (income-statement figures), columns D:M.
# ---------------------------------------------------------------
$rowData = @{
    11 = @(1617629, 2229241, 846825, 1537387, 2104278, 3058239, 1122218, 1957393, 3104953, 4208226)
    12 = @(-848514, -1217079, -459422, -806422, -1187689, -1771311, -707841, -1321093, -2007142, -2861309)
    13 = @(769115, 1012162, 387403, 730965, 916589, 1286928, 414377, 636300, 1097811, 1346917)
    14 = @(-80761, -106716, -38489, -77600, -113826, -156399, -56888, -106564, -198849, -263559)
    15 = @(0, 2414, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(240, 34, 5116, 584, 48958, 0, 25, -5188, -4434, -3677)
    17 = @(688594, 907894, 354030, 653949, 851721, 1130529, 357514, 524548, 894528, 1079681)
    18 = @(-2037, -2832, 0, 0, 0, 0, 0, -835, -3339, -5843)
    19 = @(65119, 85110, 7661, 18058, 37029, -8980, 1661, 11408, 9158, 15321)
    20 = @(751676, 990172, 361691, 672007, 888750, 1121549, 359175, 535121, 900347, 1089159)
    21 = @(-137311, -126423, -70806, -130770, -170344, -182177, -53876, 98773, -7932, -53306)
    22 = @(614365, 863749, 290885, 541237, 718406, 939372, 305299, 633894, 892415, 1035853)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(614365, 863749, 290885, 541237, 718406, 939372, 305299, 633894, 892415, 1035853)
    25 = @(1365, 1074, 362, 673, 894, 470, 380, 317, 446, 518)
    26 = @(450000, 804000, 804000, 804000, 804000, 2000000, 804000, 2000000, 2000000, 2000000)
    27 = @(307, 432, 145, 271, 359, 470, 153, 317, 446, 518)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}
